# This script applies a refreshed snapshot of the cryptocurrency price/
# volume table (columns D "Price" and E "Volume(1h)"), as produced by the
# periodic "Updated cryptos list" GitHub Actions job. Two coins (rows 16-17)
# also swapped rank position: WrappedEther now ranks above ShibaInu.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell -> new value updates to apply.
$updates = [ordered]@{
    'D2' = '61.282.41'
    'E2' = '  +0.86%  '
    'D3' = '3.430.73'
    'E3' = '  +1.82%  '
    'E4' = '  -0.01%  '
    'D5' = '576.39'
    'E5' = '  +1.41%  '
    'D6' = '144.46'
    'E6' = '  +6.41%  '
    'D7' = '3.432.87'
    'E7' = '  +1.93%  '
    'E8' = '  +0.06%  '
    'E9' = '  +1.79%  '
    'D10' = '7.61'
    'E10' = '  +0.60%  '
    'E11' = '  +3.26%  '
    'E12' = '  +2.02%  '
    'D13' = '4.024.99'
    'E13' = '  +1.88%  '
    'D14' = '28.00'
    'E14' = '  +8.14%  '
    'E15' = '  -0.89%  '
    'B16' = 'WrappedEther'
    'C16' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D16' = '3.484.82'
    'E16' = '  +3.32%  '
    'B17' = 'ShibaInu'
    'C17' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D17' = '0.0000172'
    'E17' = '  +1.75%  '
    'D18' = '61.457.00'
    'E18' = '  +0.93%  '
    'E19' = '  +8.61%  '
    'D20' = '14.16'
    'E20' = '  +3.80%  '
    'D21' = '9.48'
    'E21' = '  +3.00%  '
    'D22' = '394.11'
    'E22' = '  +6.31%  '
    'E23' = '  +3.15%  '
    'D24' = '72.88'
    'E24' = '  +3.17%  '
    'D25' = '0.996'
    'E25' = '  -0.54%  '
    'E26' = '  -0.19%  '
    'E27' = '  -0.06%  '
    'D28' = '3.585.31'
    'E28' = '  +2.21%  '
    'D29' = '0.178'
    'E29' = '  +0.67%  '
    'D30' = '7.57'
    'E30' = '  +3.83%  '
    'E31' = '  +0.32%  '
    'E32' = '  -7.87%  '
    'D33' = '8.13'
    'E33' = '  +1.87%  '
    'E34' = '  +2.31%  '
    'E35' = '  +0.01%  '
    'D36' = '23.96'
    'E36' = '  +3.20%  '
    'D37' = '3.468.30'
    'E37' = '  +2.24%  '
    'E38' = '  +4.10%  '
    'D39' = '5.09'
    'E39' = '  +0.20%  '
    'E40' = '  +1.63%  '
    'D41' = '167.35'
    'E41' = '  +1.73%  '
    'D42' = '0.0778'
    'E42' = '  +3.26%  '
    'D43' = '27.34'
    'E43' = '  +9.11%  '
    'D44' = '0.800'
    'E44' = '  +3.87%  '
    'E45' = '  +0.00%  '
    'D46' = '1.72'
    'E46' = '  +0.90%  '
    'D47' = '4.47'
    'E47' = '  +3.88%  '
    'D48' = '42.03'
    'E48' = '  +0.73%  '
    'D49' = '2.591.82'
    'E49' = '  +2.31%  '
    'E50' = '  +0.03%  '
    'D51' = '6.90'
    'E51' = '  +2.61%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Columns D (Price) and E (Volume) must stay plain text: some values
    # (e.g. "28.00", "0.800", "0.0000172") would otherwise be silently
    # reinterpreted as numbers and lose their original formatting.
    if ($cellRef -match "^[DE]") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$cellRef]
}
